$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "shift" column header
$ws.Range("H1").Value = "shift"

# Update capacity column (E) values from 10 to 3 for rows 2-5
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3

# Add shift values in new column H for rows 2-5
$ws.Range("H2").Value = 150
$ws.Range("H3").Value = 150
$ws.Range("H4").Value = 150
$ws.Range("H5").Value = 150

# Update selected cell to D7 as per diff
$ws.Range("D7").Select()
